# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.2881169905109251;  C = 0.00006708468553440206; D = 18.71679738969934;  E = 13.86384647080068;  G = 32.86882793569648 }
    3 = @{ B = 1.445647641019636;   C = 1.626987699542094;     D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    4 = @{ B = 3.272327238179451;   C = 1.626987699542094;     D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    5 = @{ B = 1.445647641019636;   C = 1.626987699542094;     D = 3.223369029078222;  E = 13.86384647080068;  G = 20.15985084044064 }
    6 = @{ B = 3.272327238179451;   C = 1.626987699542094;     D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    7 = @{ B = 0.6545652718822623;  C = 1.626987699542094;     D = 3.223369029078222;  E = 0.5333859586016987; G = 6.038307959104277 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    $ws.Range("B$row").Value = $cols.B
    $ws.Range("C$row").Value = $cols.C
    $ws.Range("D$row").Value = $cols.D
    $ws.Range("E$row").Value = $cols.E
    $ws.Range("G$row").Value = $cols.G
}
